$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the short code/function names in column B with descriptive titles
$ws.Range("B3").Value  = "Descriptions of initial datasets"
$ws.Range("B4").Value  = "Getting environmental values for the study sites"
$ws.Range("B5").Value  = "Plotting our datasets in G-space and E-space"
$ws.Range("B6").Value  = "Mahalanobis model"
$ws.Range("B7").Value  = "Plotting an ellipsoid in E-space"
$ws.Range("B8").Value  = "How to simulate a random sample of environmental conditions inside an region"
$ws.Range("B9").Value  = "How to identify points inside and outside the fundamental niche"
$ws.Range("B10").Value = "Weighted-normal model"
$ws.Range("B13").Value = "Evaluation method: accumulation curve of occurrences"

# Widen column B to fit the new, longer text
$ws.Columns.Item(2).ColumnWidth = 33

# Move the active selection to B13
[void]$ws.Range("B13").Select()
